# Regenerate sval data to filter save games.
# Directly overwrite the numeric columns B:G for rows 2-11 with the
# newly-computed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 1, 5.488907176552729)
    3  = @(1.505614041169197, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 1, 6.741336633845642)
    4  = @(0.7287194209349384, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 1, 3.034748368925986)
    5  = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538)
    6  = @(0.7287194209349384, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 0, 3.594575437922795)
    7  = @(3.182878228561681, 1.65323645889881, 16.98373111632243, 0.4998867070740569, 1, 22.31973251085698)
    8  = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 1, 5.488907176552729)
    9  = @(1.505614041169197, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 1, 3.811642989160245)
    10 = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 0, 8.418600821238126)
    11 = @(0.7287194209349384, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 0, 3.034748368925986)
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
